$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14, Q1=15, matching style of existing header cells ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows (2-25): swap I<->K and M<->O, then add new P and Q columns = 2 ---
for ($row = 2; $row -le 25; $row++) {
    $iVal = $ws.Cells.Item($row, 9).Value()
    $kVal = $ws.Cells.Item($row, 11).Value()
    $ws.Cells.Item($row, 9).Value = $kVal
    $ws.Cells.Item($row, 11).Value = $iVal

    $mVal = $ws.Cells.Item($row, 13).Value()
    $oVal = $ws.Cells.Item($row, 15).Value()
    $ws.Cells.Item($row, 13).Value = $oVal
    $ws.Cells.Item($row, 15).Value = $mVal

    $ws.Cells.Item($row, 16).Value = 2
    $ws.Cells.Item($row, 17).Value = 2
}

Write-Output "done"
